# "Colocando header nos graficos"
# Adds a header label to column A on each scenario sheet, fixes accented
# Portuguese text that had lost its diacritics, removes the obsolete
# "Teto" row from the emissions sheet, and refreshes the cost sheet header
# and figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheets 1-4 share the exact same row layout (Fonte/Tecnologia table).
# ---------------------------------------------------------------------
$rowLabels = @{
    2  = "Hidro"
    3  = "Gás Natural"
    4  = "Carvão"
    5  = "Nuclear"
    6  = "Óleos Comb"
    7  = "Biomassa"
    8  = "Eólica"
    9  = "Solar"
    10 = "Outros"
    11 = "Pot. Compl."
    12 = "GD"
}

foreach ($sheetIdx in 1..4) {
    $ws = $wb.Worksheets.Item($sheetIdx)

    # New header cell A1, styled like the rest of row 1.
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial(-4122)
    $ws.Range("A1").Value = "Fonte/Tecnologia"

    # Re-write each technology label (fixing accents) and drop the bold
    # bordered style those cells used to carry.
    foreach ($r in $rowLabels.Keys) {
        $cell = $ws.Range("A$r")
        $cell.Value = $rowLabels[$r]
        $cell.ClearFormats()
    }
}

# ---------------------------------------------------------------------
# Sheet 5 - Emissoes Totais (MtCO2eq)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("B1").Copy()
$ws5.Range("A1").PasteSpecial(-4122)
$ws5.Range("A1").Value = "Período"

$a2 = $ws5.Range("A2")
$a2.Value = "P.Médio"
$a2.ClearFormats()

$a3 = $ws5.Range("A3")
$a3.Value = "P.Crítico"
$a3.ClearFormats()

# Row 4 ("Teto") is no longer part of the table.
$ws5.Rows.Item(4).Delete()

# ---------------------------------------------------------------------
# Sheet 6 - Custo Total (bilhões de R$)
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

$ws6.Range("B1").Copy()
$ws6.Range("A1").PasteSpecial(-4122)
$ws6.Range("A1").Value = "Tipo Expansão"

# B1 becomes the text "2015" (same as the other sheets' year headers),
# not a number - copy a text "2015" cell over with values-only paste so
# the existing style/format of B1 is preserved.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B1").Copy()
$ws6.Range("B1").PasteSpecial(-4163)

$a2 = $ws6.Range("A2")
$a2.Value = "Expansão Centralizada"
$a2.ClearFormats()
$ws6.Range("B2").Value = 587

$a3 = $ws6.Range("A3")
$a3.Value = "Expansão por GD"
$a3.ClearFormats()
$ws6.Range("B3").Value = 99
